$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.881.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.661.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.26%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -1.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.107'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("E10").Value = '  -0.67%  '

$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.143.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.803.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '

$ws.Range("E16").Value = '  -0.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.650.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '341.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '

$ws.Range("E20").Value = '  -1.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.34%  '

$ws.Range("E24").Value = '  +11.17%  '

$ws.Range("E25").Value = '  +3.42%  '

$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '551.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +17.42%  '

$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("E29").Value = '  +0.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.69%  '

$ws.Range("E32").Value = '  +9.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0816'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '175.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.400'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("E39").Value = '  +2.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '172.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.18%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("E43").Value = '  -1.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.25%  '

$ws.Range("E45").Value = '  -1.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0546'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0961'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.30%  '

$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("E49").Value = '  +0.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '

$ws.Range("E51").Value = '  -0.72%  '
